# Update "想去人数" (want-to-go count) figures for several events,
# reflecting refreshed data generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 637
$wsExhibition.Range("F3").Value = 484
$wsExhibition.Range("F8").Value = 1371
$wsExhibition.Range("F9").Value = 4003

# Sheet "演出" (Performance)
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F2").Value = 56

# Sheet "全部类型" (All Types) - aggregated view, mirrors the same updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 637
$wsAll.Range("F3").Value = 484
$wsAll.Range("F8").Value = 1371
$wsAll.Range("F9").Value = 4003
$wsAll.Range("F11").Value = 56
